# Experiment order generation script
# Regenerates the per-task stimulus order sheets and rotates their names.
$wb = $excel.ActiveWorkbook

# --- Sheet 1 (was GNG_TO) becomes vSAT_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "vSAT_TO-16515889742837493"
$ws1.Range("B2").Value = "vSAT_stims-16515889742677402.csv"
$ws1.Range("B3").Value = "SAT_stims-16515889741997354.csv"
$ws1.Range("B4").Value = "vSAT_stims-16515889742517364.csv"
$ws1.Range("B5").Value = "SAT_stims-16515889742207344.csv"

# --- Sheet 2 (was NB_TO) becomes GNG_TO, shrinks from 8 rows to 4 rows ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "GNG_TO-16515889743350277"
$ws2.Range("B2").Value = "go_stims-16515889742907357.csv"
$ws2.Range("B3").Value = "GNG_stims-1651588974315739.csv"
$ws2.Range("B4").Value = "go_stims-16515889743177514.csv"
$ws2.Range("B5").Value = "GNG_stims-1651588974334025.csv"
$ws2.Rows("6:10").Delete() | Out-Null

# --- Sheet 3 (RS_TO) keeps its name pattern, new suffix, values swapped ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16515889743370278"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL_TO) keeps its name pattern, new suffix, refreshed stim files ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16515889743985598"
$ws4.Range("B2").Value = "MM_stims-16515889743650274.csv"
$ws4.Range("B3").Value = "ZM_stims-16515889743420284.csv"
$ws4.Range("B4").Value = "MM_stims-16515889743810246.csv"
$ws4.Range("B5").Value = "ZM_stims-1651588974367026.csv"
$ws4.Range("B6").Value = "MM_stims-16515889743965578.csv"
$ws4.Range("B7").Value = "ZM_stims-1651588974382025.csv"

# --- Sheet 5 (was vSAT_TO) becomes NB_TO, grows from 4 rows to 8 rows ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "NB_TO-16515889762171862"
$ws5.Range("B2").Value = "OB-16515889755591323.csv"
$ws5.Range("B3").Value = "ZB-match_3-1651588974732189.csv"
$ws5.Range("B4").Value = "ZB-match_6-16515889750929394.csv"
$ws5.Range("B5").Value = "TB-16515889762015626.csv"

$ws5.Range("A2:B2").Copy() | Out-Null
$ws5.Range("A6:B10").PasteSpecial(-4122) | Out-Null

$ws5.Range("A6").Value = 4
$ws5.Range("B6").Value = "ZB-match_9-16515889750535989.csv"
$ws5.Range("A7").Value = 5
$ws5.Range("B7").Value = "TB-16515889760729644.csv"
$ws5.Range("A8").Value = 6
$ws5.Range("B8").Value = "OB-16515889752536712.csv"
$ws5.Range("A9").Value = 7
$ws5.Range("B9").Value = "TB-16515889761036978.csv"
$ws5.Range("A10").Value = 8
$ws5.Range("B10").Value = "OB-16515889753912845.csv"
